# Replace the Google Drive hyperlink with a new URL, stored as plain
# (non-hyperlinked) text, matching the target edit:
#   - old hyperlink (field + "Hyperlink" styled run) pointing at
#     .../135MVw6PmlBONxGimvYze47yCH9tvPomd/... is removed
#   - replaced by a plain run containing the new URL
#     .../1ZU7I_AmjiUhPLC6i8OVpwefCWKXLNTLX/...

$d = $word.ActiveDocument

$oldUrl = "https://drive.google.com/file/d/135MVw6PmlBONxGimvYze47yCH9tvPomd/view?usp=sharing"
$newUrl = "https://drive.google.com/file/d/1ZU7I_AmjiUhPLC6i8OVpwefCWKXLNTLX/view?usp=sharing"

# Walk the hyperlinks backwards (so deleting doesn't disturb not-yet-visited
# indices) and replace any one whose address/display text is the old URL.
for ($i = $d.Hyperlinks.Count; $i -ge 1; $i--) {
    $h = $d.Hyperlinks.Item($i)

    if ($h.Address -eq $oldUrl -or $h.TextToDisplay -eq $oldUrl) {
        $start = $h.Range.Start
        $end = $h.Range.End

        # Remove the hyperlinked run (this deletes the hyperlink field
        # together with its "Hyperlink" styled text run).
        $killRange = $d.Range($start, $end)
        $killRange.Delete()

        # Insert the new URL as ordinary, unstyled text at the same spot.
        $insertRange = $d.Range($start, $start)
        $insertRange.InsertAfter($newUrl)
    }
}
